$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 06:55"

# 2. Row 6 - India (rank 10)
$ws.Range("B6").Value = 3310234
$ws.Range("C6").Value = 2485
$ws.Range("D6").Value = 2523771
$ws.Range("E6").Value = 725834

# 3. Row 60 - Venezuela (rank 64)
$ws.Range("B60").Value = 41965
$ws.Range("D60").Value = 32931
$ws.Range("E60").Value = 8683
$ws.Range("H60").Value = 351

# 4. Row 62 - Uzbekistan (rank 66)
$ws.Range("B62").Value = 40054
$ws.Range("C62").Value = 90
$ws.Range("D62").Value = 36562
$ws.Range("E62").Value = 3199
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 293

# 5. Row 123 - Tailandia (rank 127)
$ws.Range("B123").Value = 3404
$ws.Range("C123").Value = 1
$ws.Range("E123").Value = 109

# 6. Rows 143-146: new "Jamaica" entry inserted right after Bahamas, cascading
#    Aruba / Jordania / Malta down one row (the old, lower "Jamaica" row is removed).
# Row 143 becomes Jamaica, with fresh data
$ws.Range("A143").Value = "Jamaica"
$ws.Range("B143").Value = 1804
$ws.Range("C143").Value = 72
$ws.Range("D143").Value = 846
$ws.Range("E143").Value = 939
$ws.Range("H143").Value = 19

# Row 144 becomes Aruba, carrying the old row-143 data
$ws.Range("A144").Value = "Aruba"
$ws.Range("B144").Value = 1760
$ws.Range("D144").Value = 587
$ws.Range("E144").Value = 1165
$ws.Range("H144").Value = 8

# Row 145 becomes Jordania, carrying the old row-144 data
$ws.Range("A145").Value = "Jordania"
$ws.Range("B145").Value = 1756
$ws.Range("D145").Value = 1355
$ws.Range("E145").Value = 386
$ws.Range("H145").Value = 15

# Row 146 becomes Malta, carrying the old row-145 data
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 1751
$ws.Range("D146").Value = 1077
$ws.Range("E146").Value = 664
$ws.Range("H146").Value = 10

# 7. Row 170 - Birmania (rank 174)
$ws.Range("B170").Value = 586
$ws.Range("C170").Value = 6
$ws.Range("E170").Value = 235

# 8. Row 183 - Mongolia (rank 187)
$ws.Range("B183").Value = 301
$ws.Range("C183").Value = 1
$ws.Range("D183").Value = 291
$ws.Range("E183").Value = 10

# 9. Row 197 - Curazao (rank 201)
$ws.Range("B197").Value = 53
$ws.Range("C197").Value = 4
$ws.Range("E197").Value = 18

Write-Host "Edits applied"
